$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.277.88'
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").Value = '2.065.53'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.620'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.85%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.79'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.383'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0763'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.27%  '
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("D12").Value = '2.369.86'
$ws.Range("E12").Value = '  -0.68%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.63'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.81'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.777'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.15'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("D17").Value = '2.065.01'
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("D18").Value = '37.222.15'
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.37'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.46%  '
$ws.Range("D21").Value = '0.0₃0813'
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.43'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.40'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.42'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.43'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.87%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.65%  '
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.126'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.118'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.47'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.20%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0616'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.21%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.60'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.20%  '
$ws.Range("E35").Value = '  -4.50%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.77'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.32%  '
$ws.Range("E38").Value = '  -3.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.65'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.31%  '
$ws.Range("E40").Value = '  -0.92%  '
$ws.Range("D41").Value = '1.478.00'
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.05'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.55%  '
$ws.Range("E43").Value = '  +0.36%  '
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0933'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.17%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.17'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.74%  '
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("E47").Value = '  -0.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.81%  '
$ws.Range("E49").Value = '  +1.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.16'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.98'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.30%  '
